$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "invalid number of offers" row (row 4) - mirrors rows 1-3
$ws.Range("A4").Value = "ABCHJUH"
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = "909ikokujyhtgt*"
$ws.Range("D4").Value = "JKJKUHY/////\\\\\%^%gyvb"
$ws.Range("E4").Value = "iojkjkjhjhjhjjhj"
$ws.Range("F4").Value = "ioiojkhjghfgfgghbn"
$ws.Range("G4").Value = "jkjkhjhh)))))"
$ws.Range("H4").Value = "(((hjnmnmnmm####"
$ws.Range("I4").Value = "hjhj.uiuiuisdksd"

# New column (H) needs an explicit width, same as the other bestFit columns
$ws.Columns("H").ColumnWidth = 16.666

# Final selection left on column G (whole column) by the author
$ws.Columns("G").Select() | Out-Null
